$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "tool_description"
$ws.Range("C1").Value = "contact_name"
$ws.Range("D1").Value = "contact_email"
$ws.Range("E1").Value = "is_maintained"
$ws.Range("F1").Value = "relevant_diseases"
$ws.Range("K1").Value = "required_expertise"
$ws.Range("L1").Value = "tool_type"
$ws.Range("M1").Value = "input_type"
$ws.Range("N1").Value = "docs_link"
$ws.Range("O1").Value = "source_link"
$ws.Range("Q1").Value = "github_link"
$ws.Range("R1").Value = "is_complete"
$ws.Range("S1").Value = "pkg_dev_assessment"
$ws.Range("T1").Value = "overall_assessment"
